$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("A18").Value = "WGG 02"
$ws.Range("B18").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C18").Value = "20-01-2026"
$ws.Range("D18").Value = 286962
$ws.Range("E18").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F18").Value = 34400000000
$ws.Range("G18").Value = "NEFT"
$ws.Range("H18").Value = "SBIN0003229"
$ws.Range("I18").Value = "AAAFW8862C"
$ws.Range("J18").Value = "32AAAFW8862C1Z9"
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = "252f7cfc-b58e-4193-a6b3-455a0d3a265d"
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = ""
$ws.Range("O18").Value = ""
$ws.Range("P18").Value = ""
$ws.Range("Q18").Value = ""
$ws.Range("R18").Value = ""
$ws.Range("S18").Value = ""
$ws.Range("T18").Value = ""
$ws.Range("U18").Value = "pending"
$ws.Range("V18").Value = 500
$ws.Range("W18").Value = ""
$ws.Range("X18").Value = "PAYMENT TESTING RPA_UNIQUE_ID : a339badc-3f04-410f-890d-07c3c5d16a78"
$ws.Range("Y18").Value = "HO"
$ws.Range("Z18").Value = 0
$ws.Range("AA18").Value = "midhuncraju12@gmail.com"
$ws.Range("AB18").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC18").Value = 0
$ws.Range("AD18").Value = 0
$ws.Range("AE18").Value = 0
$ws.Range("AF18").Value = ""
$ws.Range("AG18").Value = ""
$ws.Range("AH18").Value = ""
$ws.Range("AI18").Value = ""
$ws.Range("AJ18").Value = ""
$ws.Range("AK18").Value = ""
$ws.Range("AL18").Value = ""
$ws.Range("AM18").Value = ""
$ws.Range("AN18").Value = ""
$ws.Range("AO18").Value = ""

# Row 19
$ws.Range("A19").Value = "WGG 02"
$ws.Range("B19").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C19").Value = "20-01-2026"
$ws.Range("D19").Value = 286962
$ws.Range("E19").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F19").Value = 34400000000
$ws.Range("G19").Value = "NEFT"
$ws.Range("H19").Value = "SBIN0003229"
$ws.Range("I19").Value = "AAAFW8862C"
$ws.Range("J19").Value = "32AAAFW8862C1Z9"
$ws.Range("K19").Value = ""
$ws.Range("L19").Value = "56fac636-12e0-4b8d-9d8f-d5a2596dd10e"
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = ""
$ws.Range("O19").Value = ""
$ws.Range("P19").Value = ""
$ws.Range("Q19").Value = ""
$ws.Range("R19").Value = ""
$ws.Range("S19").Value = ""
$ws.Range("T19").Value = ""
$ws.Range("U19").Value = "pending"
$ws.Range("V19").Value = 1500
$ws.Range("W19").Value = ""
$ws.Range("X19").Value = "PAYMENT TESTING RPA_UNIQUE_ID : 4e2705d4-7123-4dc9-b863-668f4586341c"
$ws.Range("Y19").Value = "HO"
$ws.Range("Z19").Value = 0
$ws.Range("AA19").Value = "midhuncraju12@gmail.com"
$ws.Range("AB19").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC19").Value = 0
$ws.Range("AD19").Value = 0
$ws.Range("AE19").Value = 0
$ws.Range("AF19").Value = ""
$ws.Range("AG19").Value = ""
$ws.Range("AH19").Value = ""
$ws.Range("AI19").Value = ""
$ws.Range("AJ19").Value = ""
$ws.Range("AK19").Value = ""
$ws.Range("AL19").Value = ""
$ws.Range("AM19").Value = ""
$ws.Range("AN19").Value = ""
$ws.Range("AO19").Value = ""
